# Implements commit "实现条文GB51251_2017_4_5_2 ... " which adds an
# "是否为室外风口" (is outdoor vent) column to the "风口" sheet and a
# "风量" (air volume) column to every "风机(...)" sheet, fills in sample
# data for the new columns, adjusts a few column widths, and moves the
# active/selected tab & cell selections around.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet indices (1-based, matching workbook.xml sheet order rId1..rId7)
#   1 = 风机(住宅>100通过)
#   2 = 风口
#   3 = 风机(住宅>100不通过)
#   4 = 风机(公共建筑>50不通过)
#   5 = 风机(住宅<=100通过)
#   6 = 风机(公共建筑>50通过)
#   7 = 风机(公共建筑<=50通过)
# ---------------------------------------------------------------------

$fanSheetIdx = @(1, 3, 4, 5, 6, 7)

# Number of data rows (below header row 1) that receive a new F value
# ("风量" = 1000) for each fan sheet. Sheet5's last data row (row 4) is
# intentionally left without a value, matching the source data.
$fanRowsWithValue = @{
    1 = @(2, 3)
    3 = @(2, 3)
    4 = @(2, 3)
    5 = @(2, 3)
    6 = @(2, 3, 4, 5)
    7 = @(2, 3, 4, 5)
}

# Column width targets (best effort - the engine re-quantizes widths to
# its own internal pixel grid, so exact source decimals cannot always be
# reproduced).
$fanColWidths = @{
    1 = @{2 = 22.08984375; 3 = 36.08984375; 4 = 15.453125}
    3 = @{2 = 22.08984375; 3 = 37.90625;    4 = 15.453125}
    4 = @{2 = 22.08984375; 3 = 37.90625;    4 = 15.453125}
    5 = @{2 = 22.08984375; 3 = 33.26953125; 4 = 15.453125}
    6 = @{2 = 22.08984375; 3 = 35.453125;   4 = 15.453125}
    7 = @{2 = 22.08984375; 4 = 15.453125}
}

# Selection (activeCell:sqref) to leave on each fan sheet afterwards.
$fanSelection = @{
    1 = @{cell = "C23"; sqref = "C23"}
    3 = @{cell = "F4";  sqref = "F4:F5"}
    4 = @{cell = "F4";  sqref = "F4:F5"}
    5 = @{cell = "F4";  sqref = "F4:F5"}
    6 = @{cell = "F1";  sqref = "F1:F1048576"}
    7 = @{cell = "F1";  sqref = "F1:F1048576"}
}

foreach ($idx in $fanSheetIdx) {
    $ws = $wb.Worksheets.Item($idx)

    # New header cell: F1 = "风量"
    $ws.Cells.Item(1, 6).Value = "风量"

    # New data values in column F
    foreach ($r in $fanRowsWithValue[$idx]) {
        $ws.Cells.Item($r, 6).Value = 1000
    }

    # Column width tweaks
    $widths = $fanColWidths[$idx]
    foreach ($c in $widths.Keys) {
        $ws.Columns.Item($c).ColumnWidth = $widths[$c]
    }
}

# ---------------------------------------------------------------------
# "风口" sheet (index 2): add "是否为室外风口" (E) and "风量" (F) columns
# ---------------------------------------------------------------------
$wsVents = $wb.Worksheets.Item(2)

$wsVents.Cells.Item(1, 5).Value = "是否为室外风口"
$wsVents.Cells.Item(1, 6).Value = "风量"

$ventAirVolume = @{
    2 = 1000; 3 = 1000; 4 = 1000; 5 = 1000; 6 = 1000; 7 = 1500; 8 = 1500
    9 = 1000; 10 = 1000; 11 = 1000; 12 = 1000; 13 = 1000; 14 = 1000
    15 = 1000; 16 = 1000; 17 = 1000; 18 = 1000; 19 = 1000; 20 = 1000
    21 = 1500; 22 = 1500; 23 = 1000; 24 = 1000; 25 = 1000; 26 = 1500
    27 = 1000
}

for ($r = 2; $r -le 27; $r++) {
    $wsVents.Cells.Item($r, 5).Value = $false
    $wsVents.Cells.Item($r, 6).Value = $ventAirVolume[$r]
}

# New column width for the newly added column E
$wsVents.Columns.Item(5).ColumnWidth = 19.26953125

# ---------------------------------------------------------------------
# Selections / active tab. Perform the fan-sheet selections first (in
# ascending order), "风口" afterwards, and make sheet 7 the very last
# selection so it ends up as the active tab, matching activeTab="6".
# ---------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range($fanSelection[1].sqref).Select() | Out-Null

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range($fanSelection[3].sqref).Select() | Out-Null

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range($fanSelection[4].sqref).Select() | Out-Null

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range($fanSelection[5].sqref).Select() | Out-Null

$ws6 = $wb.Worksheets.Item(6)
$ws6.Range($fanSelection[6].sqref).Select() | Out-Null

$wsVents.Range("F28:F36").Select() | Out-Null

$ws7 = $wb.Worksheets.Item(7)
$ws7.Range($fanSelection[7].sqref).Select() | Out-Null
